# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the formatting of the existing header row (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for the Save column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
